# Actualización automática 2025-07-01 11:50:08
#
# Inserts a new advisor "ANGULO PARRALES CARMEN" as the first data row
# (row 2) in both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets,
# pushing the existing advisor rows (and the trailing summary row) down
# by one, and updates the summary row's "0 de 6" / "0 de 7" wording
# where present.
#
# Rows are shifted with Range.Copy(Destination) (source row -> row+1),
# walking bottom-up so a row is always read before it gets overwritten.
# Copy carries the per-cell style along with the value, so every moved
# cell keeps exactly the formatting it already had - including the
# summary row, which lands one row below the sheet's previous used
# range and would otherwise pick up no style at all.
#
# (Rows.Insert() was avoided: in this host it always stamps the new row
# with the format of the row above - here the bold/centered header -
# which is wrong for a plain data row.)

$wb = $excel.ActiveWorkbook

function Shift-AdvisorSheet {
    param($SheetName, $LastDataCol, $SummaryIsText, $OldAdvisorCount, $NewAdvisorCount)

    $ws = $wb.Worksheets.Item($SheetName)
    $lastColLetter = [char](64 + $LastDataCol)

    $firstDataRow = 2
    $lastAdvisorRow = $firstDataRow + $OldAdvisorCount - 1   # old last advisor row (e.g. 7)
    $oldSummaryRow = $lastAdvisorRow + 1                      # old summary row (e.g. 8)
    $newSummaryRow = $firstDataRow + $NewAdvisorCount         # new summary row (e.g. 9)

    # Shift every existing row (advisors + summary) down by one, working
    # from the bottom up so a row is copied before it is overwritten.
    # The summary row only has cells from column C onward (no ASESOR /
    # CLIENTE columns), so it is copied as C:<last> only - otherwise the
    # copy would materialise blank A/B cells that shouldn't exist.
    for ($r = $oldSummaryRow; $r -ge $firstDataRow; $r--) {
        $destRow = $r + 1
        if ($r -eq $oldSummaryRow) {
            $srcRange = $ws.Range("C" + $r + ":" + $lastColLetter + $r)
            $dstRange = $ws.Range("C" + $destRow + ":" + $lastColLetter + $destRow)
        } else {
            $srcRange = $ws.Range("A" + $r + ":" + $lastColLetter + $r)
            $dstRange = $ws.Range("A" + $destRow + ":" + $lastColLetter + $destRow)
        }
        $srcRange.Copy($dstRange)
    }

    # Write the new first advisor row (its cells already carry the
    # correct "data row" style from the template, untouched above).
    $ws.Cells.Item($firstDataRow, 1).Value = "OFICINA-CATAECSA"
    $ws.Cells.Item($firstDataRow, 2).Value = "ANGULO PARRALES CARMEN"
    for ($col = 3; $col -le $LastDataCol; $col++) {
        $ws.Cells.Item($firstDataRow, $col).Value = 0
    }

    # Update the "0 de N" wording on the relocated summary row, if present.
    if ($SummaryIsText) {
        $oldText = "0 de " + $OldAdvisorCount
        $newText = "0 de " + $NewAdvisorCount
        for ($col = 3; $col -le $LastDataCol; $col++) {
            $ws.Cells.Item($newSummaryRow, $col).Value = $newText
        }
    }
}

# "VENTAS POR GRUPO": columns C..R (18 cols total), summary row holds "0 de 6" -> "0 de 7" text.
Shift-AdvisorSheet "VENTAS POR GRUPO" 18 $true 6 7

# "VENTA MENSUAL": columns C..G (7 cols total), summary row holds numeric totals (unchanged, just relocated).
Shift-AdvisorSheet "VENTA MENSUAL" 7 $false 6 7
